# Update the "Riders" (C) and "Average" (D) columns on the Ridership sheet
# with the new Madigan bike hours figures. The chart on this sheet plots
# these ranges directly, so its cached values will be refreshed when the
# workbook recalculates / is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

$newRiders = @(91,158,209,247,266,136,86,265,222,218,241,289,108,85,238,177,275,246,176,125,211,176,204,220,236,111,131,286,241)
$newAverage = @(229.3,208.6,193.5,217.78,231.3,118.27,103.33,232.55,209.82,195.73,220.1,236.55,117.42,101.5,212.17,194.17,225.09,237.33,121.92,103.64,230.75,209.38,194.92,224.67,237.23,121.14,105.92,235,211.64)

for ($i = 0; $i -lt $newRiders.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newRiders[$i]
    $ws.Cells.Item($row, 4).Value = $newAverage[$i]
}

$excel.CalculateFullRebuild()
